$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.901.84"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "'1.637.55"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D5").Value = "'213.53"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'23.54"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "'0.0614"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").Value = "'0.0874"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "'1.870.12"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "'1.633.34"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.10"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.574"
$ws.Range("E15").Value = "  +3.91%  "
$ws.Range("D16").Value = "'65.94"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "'27.895.59"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").Value = "'231.23"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").Value = "'7.61"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'10.78"
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = "  -3.25%  "
$ws.Range("D25").Value = "'151.63"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").Value = "'6.92"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").Value = "'15.70"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'3.33"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("D33").Value = "'3.09"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("D34").Value = "'1.401.33"
$ws.Range("E34").Value = "  -5.00%  "
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "'0.917"
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("E43").Value = "  +4.89%  "
$ws.Range("D44").Value = "'66.12"
$ws.Range("E44").Value = "  -2.72%  "
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "'1.778.80"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").Value = "'87.93"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").Value = "'7.63"
$ws.Range("E51").Value = "  -1.13%  "
